$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

# Row 2
$ws.Range("D2").Value = "63.833.45"
$ws.Range("E2").Value = "  +1.15%  "

# Row 3
$ws.Range("D3").Value = "3.326.56"
$ws.Range("E3").Value = "  +2.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.81%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "3.325.69"
$ws.Range("E8").Value = "  +2.47%  "

# Row 9
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
$ws.Range("E10").Value = "  +1.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.57%  "

# Row 12
$ws.Range("E12").Value = "  +1.58%  "

# Row 13
$ws.Range("E13").Value = "  +0.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "

# Row 15
$ws.Range("D15").Value = "3.873.28"
$ws.Range("E15").Value = "  +2.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17
$ws.Range("D17").Value = "3.325.97"
$ws.Range("E17").Value = "  +2.43%  "

# Row 18
$ws.Range("D18").Value = "63.895.05"
$ws.Range("E18").Value = "  +1.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "

# Row 21
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("E22").Value = "  +2.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.65%  "

# Row 29
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "

# Row 31
$ws.Range("E31").Value = "  +2.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.53%  "

# Row 33
$ws.Range("E33").Value = "  -1.02%  "

# Row 34
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.68%  "

# Row 36
$ws.Range("E36").Value = "  +3.73%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0753"
$ws.Range("E37").Value = "  +5.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.09%  "

# Row 40
$ws.Range("D40").Value = "3.120.79"
$ws.Range("E40").Value = "  +4.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "431.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.84%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

# Row 43
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.117"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.78%  "

# Row 45
$ws.Range("E45").Value = "  +0.05%  "

# Row 46
$ws.Range("E46").Value = "  +4.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

# Row 51
$ws.Range("E51").Value = "  -0.42%  "

